$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "Trends Status" sheet - update values
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Trends Status")

$ws1.Range("B2").Value = 2
$ws1.Range("C2").Value = 20
$ws1.Range("D2").Value = 5.3
$ws1.Range("E2").Value = 26

$ws1.Range("B3").Value = 4
$ws1.Range("C3").Value = 20
$ws1.Range("D3").Value = 10.5
$ws1.Range("E3").Value = 26

$ws1.Range("B4").Value = 11
$ws1.Range("C4").Value = 25
$ws1.Range("D4").Value = 28.9
$ws1.Range("E4").Value = 32.5

$ws1.Range("B5").Value = 6
$ws1.Range("C5").Value = 4
$ws1.Range("D5").Value = 15.8
$ws1.Range("E5").Value = 5.2

$ws1.Range("B6").Value = 15
$ws1.Range("C6").Value = 8
$ws1.Range("D6").Value = 39.5
$ws1.Range("E6").Value = 10.4

$ws1.Range("B7").Value = 43
$ws1.Range("C7").Value = 86

# ---------------------------------------------------------------------------
# 2) "Priority Status" sheet - update values
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Priority Status")

$ws3.Range("B2").Value = 103
$ws3.Range("B3").Value = 286
$ws3.Range("B4").Value = 554

# ---------------------------------------------------------------------------
# 3) "Species qualification" sheet - update text + values
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Species qualification")

$ws4.Range("A2").Value = "SoIB Assessment"
$ws4.Range("B2").Value = 413

$ws4.Range("B3").Value = 81
$ws4.Range("C3").Value = 38

$ws4.Range("B4").Value = 163
$ws4.Range("C4").Value = 77

# ---------------------------------------------------------------------------
# 4) Add the new "Major update - High Priority " sheet, duplicating the
#    original "High Priority break-up" data (positioned as the last sheet
#    in the workbook), BEFORE the old sheet's own content gets overwritten
#    below. Use Add(Before, After) with After = last sheet so it lands at
#    the end directly (avoids a separate .Move() call on a freshly created
#    sheet, whose handle does not track the sheet reliably across a move
#    in this COM-interop layer).
# ---------------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("High Priority break-up")

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "Major update - High Priority "

$newSheet.Range("A1").Value = "Break-up"
$newSheet.Range("B1").Value = "High Species (no.)"
$newSheet.Range("C1").Value = "High Species (perc.)"
$newSheet.Range("D1").Value = "New High Species (no.)"
$newSheet.Range("E1").Value = "New High Species (perc.)"
$newSheet.Range("A1:E1").Font.Bold = $true
$newSheet.Range("A1:E1").HorizontalAlignment = -4108

$newSheet.Range("A2").Value = "Trend New"
$newSheet.Range("B2").Value = 25
$newSheet.Range("C2").Value = 55.6
$newSheet.Range("D2").Value = 25
$newSheet.Range("E2").Value = 55.6

$newSheet.Range("A3").Value = "IUCN"
$newSheet.Range("B3").Value = 20
$newSheet.Range("C3").Value = 44.4
$newSheet.Range("D3").Value = 20
$newSheet.Range("E3").Value = 44.4

# ---------------------------------------------------------------------------
# 5) Rename the original "High Priority break-up" sheet and overwrite its
#    data with the new "Interannual update" figures.
# ---------------------------------------------------------------------------
$ws5.Name = "Interannual update - High Pri"

$ws5.Range("A2").Value = "Trend New"
$ws5.Range("B2").Value = 68
$ws5.Range("C2").Value = 66
$ws5.Range("D2").Value = 68
$ws5.Range("E2").Value = 89.5

$ws5.Range("A3").Value = "Trend Different"
$ws5.Range("B3").Value = 8
$ws5.Range("C3").Value = 7.8
$ws5.Range("D3").ClearContents()
$ws5.Range("E3").ClearContents()

$ws5.Range("A4").Value = "IUCN"
$ws5.Range("B4").Value = 27
$ws5.Range("C4").Value = 26.2
$ws5.Range("D4").Value = 8
$ws5.Range("E4").Value = 10.5
